# Update Betfair back/lay odds for 2026-01-28 fixtures (rows 2-13 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.55
$ws.Range("G2").Value = 1.74
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 8.199999999999999
$ws.Range("J2").Value = 3.9
$ws.Range("K2").Value = 4.7
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1.93
$ws.Range("Q2").Value = 1.86
$ws.Range("G3").Value = 3.8
$ws.Range("H3").Value = 2.24
$ws.Range("F4").Value = 1.35
$ws.Range("G4").Value = 1.46
$ws.Range("H4").Value = 8.6
$ws.Range("I4").Value = 11
$ws.Range("J4").Value = 5.1
$ws.Range("K4").Value = 6.2
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 2.42
$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 2.42
$ws.Range("Q4").Value = 1.55
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 2.18
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000
$ws.Range("H5").Value = 2.34
$ws.Range("F6").Value = 2.36
$ws.Range("H6").Value = 2.9
$ws.Range("J6").Value = 2.92
$ws.Range("P6").Value = 1.55
$ws.Range("H7").Value = 9
$ws.Range("L7").Value = 1.41
$ws.Range("F10").Value = 1.52
$ws.Range("I10").Value = 17
$ws.Range("G11").Value = 2.5
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 2.98
$ws.Range("P11").Value = 1.52
$ws.Range("Q11").Value = 2.44
$ws.Range("H12").Value = 8.6
$ws.Range("K12").Value = 4.3
$ws.Range("Q12").Value = 2.3
$ws.Range("J13").Value = 4.4
$ws.Range("Q13").Value = 1.81
